# Generate Report for Handoff
# Replace the two tracked source files (UUID-based markdown names) and refresh
# their localization status/handoff metadata across all three sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "4a77df15-d595-4b47-9f73-1ebde914afab.md"
$oldFile2 = "6eb92914-7f2b-4f0d-9ef6-95f1db5d4946.md"
$newFile1 = "09ce317a-cf1b-4e00-9824-108efd02fd8d.md"
$newFile2 = "ffffaa324517-8580-4c7a-a5ec-17b83947e2b8.md"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-13 01:18:08"
$newHbDate = "0001-01-01 00:00:00"

$zhXlf = "09ce317a-cf1b-4e00-9824-108efd02fd8d.30bae8048b2f6e2a190d0f4fcee704f5df3250a8.zh-cn.xlf"
$deXlf = "09ce317a-cf1b-4e00-9824-108efd02fd8d.30bae8048b2f6e2a190d0f4fcee704f5df3250a8.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newHoDate

$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\" + $newFile1
$wsOverview.Hyperlinks.Item(2).TextToDisplay = "e2e\" + $newFile2

$wsOverview.Range("E1:F1").ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $newHoDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $newHbDate

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $newHoDate
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $newHbDate

$wsZh.Hyperlinks.Item(2).Delete()
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsZh.Hyperlinks.Item(1).Range().Value = $newFile1
$newA3Link = $wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a19c717b59edffef0159b089b9d990e9dc0dec4a/e2e/4a77df15-d595-4b47-9f73-1ebde914afab.md", "", "", $newFile2)

$wsZh.Range("C1").ColumnWidth = 17.2159881591797
$wsZh.Range("I1").ColumnWidth = 18.6506053379604
$wsZh.Range("J1").ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $newHbDate

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $newHoDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $newHbDate

$wsDe.Hyperlinks.Item(2).Delete()
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsDe.Hyperlinks.Item(1).Range().Value = $newFile1
$newA3LinkDe = $wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cd8970b58fe2fa33e7c8692c6a91b3238591a177/e2e/4a77df15-d595-4b47-9f73-1ebde914afab.md", "", "", $newFile2)

$wsDe.Range("C1").ColumnWidth = 17.2159881591797
$wsDe.Range("I1").ColumnWidth = 18.6506053379604
$wsDe.Range("J1").ColumnWidth = 21.7054770333426
